$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GestionActivos")

# Update "Pasos a seguir" for CP_GESTACT_002 (row 3, column F) with the new
# expanded step list that inserts a new step 2 and renumbers the rest.
$newSteps = '1.Clic en botón "Seleccionar entidad"' + "`n" + `
    '2.Seleccionar la fila que contiene el texto "elemento secundario' + "`n" + `
    '3.Clic en botón "Siguiente"' + "`n" + `
    '4.Seleccionar "ont"' + "`n" + `
    '5.Clic en botón "Siguiente"' + "`n" + `
    '6.Seleccionar fila con ID 9 "FAILED"' + "`n" + `
    '7.Clic en botón "FINALIZAR"'

$ws.Range("F3").Value = $newSteps
